$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.934999999999999
$ws.Range("D7").Value = -7.22289999999999
$ws.Range("C8").Value = -11.04189999999999
$ws.Range("A12").Value = -21.63070000000002
$ws.Range("C12").Value = -12.0979
$ws.Range("C14").Value = -12.2229
$ws.Range("D19").Value = -8.391799999999991
$ws.Range("E19").Value = 13.7283
$ws.Range("D21").Value = -8.014399999999991
$ws.Range("C22").Value = -10.77969999999999
$ws.Range("D24").Value = -7.826999999999996
